# Cyclic shift of taxon data among rows 11-14:
# New row 11 = old row 12 data (for columns A,B,D,E,F,G,H,Q,R)
# New row 12 = old row 13 data
# New row 13 = old row 14 data
# New row 14 = old row 11 data
# All other columns (C, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(11, 12, 13, 14)

# Capture the original values for the columns that move, keyed by row number.
# NOTE: reading the bare ".Value" property in this runtime yields a reflection
# description string instead of the cell contents, so read via ".Value()".
$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = @{
        A = $ws.Range("A$r").Value()
        B = $ws.Range("B$r").Value()
        D = $ws.Range("D$r").Value()
        E = $ws.Range("E$r").Value()
        F = $ws.Range("F$r").Value()
        G = $ws.Range("G$r").Value()
        H = $ws.Range("H$r").Value()
        Q = $ws.Range("Q$r").Value()
        R = $ws.Range("R$r").Value()
    }
}

# Mapping: destination row -> source row (data moves "up" one row, wrapping 11 -> 14)
$mapping = @{
    11 = 12
    12 = 13
    13 = 14
    14 = 11
}

foreach ($dest in $rows) {
    $src = $mapping[$dest]
    $vals = $orig[$src]

    $ws.Range("A$dest").Value = $vals.A
    $ws.Range("B$dest").Value = $vals.B
    $ws.Range("D$dest").Value = $vals.D
    $ws.Range("E$dest").Value = $vals.E
    $ws.Range("F$dest").Value = $vals.F
    $ws.Range("G$dest").Value = $vals.G
    $ws.Range("H$dest").Value = $vals.H
    $ws.Range("Q$dest").Value = $vals.Q
    $ws.Range("R$dest").Value = $vals.R
}
